# WeaponDB_Sheet: add a "_particleID" column (W) with particle ids for
# the Bow (rows 2-4) and Wand (rows 5-7) weapon rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("W1").Value = "_particleID"

# Data — Bow rows get 71000002, Wand rows get 71000001.
$ws.Range("W2").Value = 71000002
$ws.Range("W3").Value = 71000002
$ws.Range("W4").Value = 71000002
$ws.Range("W5").Value = 71000001
$ws.Range("W6").Value = 71000001
$ws.Range("W7").Value = 71000001

# New column width to match the rest of the sheet's manual sizing
# (target stored width ~19.25 chars; the host's column-width pixel
# quantization snaps the nearby setter inputs to the closest step).
$ws.Columns("W").ColumnWidth = 18.6

# Selection left where the author's cursor ended up after adding the column.
[void]$ws.Range("W3").Select()
